# Restructure the "Utlån Statistikk" sheet from a wide (horizontal) layout
# into a tall (vertical) two-column table with headers "År" / "Verdi".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old horizontal layout (row 1: years in B1:G1, row 2: label
# in A2 + values in B2:G2) before laying down the new data.
$ws.Range("A1:G2").Clear()

# New headers.
$ws.Range("A1").Value = "År"
$ws.Range("B1").Value = "Verdi"

# New vertical year/value pairs.
$years  = @(1905, 1910, 1915, 1920, 1925, 1930)
$values = @(340, 449, 656, 1732, 1805, 1491)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("B2").Select()
